# Updates the cryptocurrency price/volume snapshot values and the
# "Hora" (hour) column from 22 -> 23, matching the refreshed data pulled
# by the scheduled GitHub Actions job.
#
# Values in columns D (Price) and G (Hora) are stored as text in this
# workbook (not numbers), so each assignment below uses a leading
# apostrophe to force Excel to keep them as text, preserving formats
# such as leading/trailing zeros (e.g. "0.00001399").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'271.63"
$ws.Range("G2").Value = "'23"
$ws.Range("D3").Value = "'22.76"
$ws.Range("G3").Value = "'23"
$ws.Range("D4").Value = "'6.342"
$ws.Range("G4").Value = "'23"
$ws.Range("D5").Value = "'0.06207"
$ws.Range("G5").Value = "'23"
$ws.Range("D6").Value = "'3.651"
$ws.Range("G6").Value = "'23"
$ws.Range("D7").Value = "'6.658"
$ws.Range("G7").Value = "'23"
$ws.Range("D8").Value = "'1.383"
$ws.Range("G8").Value = "'23"
$ws.Range("D9").Value = "'0.8300"
$ws.Range("G9").Value = "'23"
$ws.Range("D10").Value = "'0.01380"
$ws.Range("G10").Value = "'23"
$ws.Range("D11").Value = "'0.1602"
$ws.Range("G11").Value = "'23"
$ws.Range("G12").Value = "'23"
$ws.Range("D13").Value = "'0.03434"
$ws.Range("G13").Value = "'23"
$ws.Range("D14").Value = "'0.03171"
$ws.Range("G14").Value = "'23"
$ws.Range("D15").Value = "'0.09341"
$ws.Range("G15").Value = "'23"
$ws.Range("D16").Value = "'3.867"
$ws.Range("G16").Value = "'23"
$ws.Range("D17").Value = "'0.001636"
$ws.Range("G17").Value = "'23"
$ws.Range("G18").Value = "'23"
$ws.Range("D19").Value = "'0.006400"
$ws.Range("G19").Value = "'23"
$ws.Range("D20").Value = "'0.005688"
$ws.Range("G20").Value = "'23"
$ws.Range("G21").Value = "'23"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("G22").Value = "'23"
$ws.Range("D23").Value = "'3.720"
$ws.Range("G23").Value = "'23"
$ws.Range("G24").Value = "'23"
$ws.Range("D25").Value = "'0.3348"
$ws.Range("G25").Value = "'23"
$ws.Range("G26").Value = "'23"
$ws.Range("D27").Value = "'0.0002703"
$ws.Range("G27").Value = "'23"
$ws.Range("G28").Value = "'23"
$ws.Range("G29").Value = "'23"
$ws.Range("G30").Value = "'23"
$ws.Range("G31").Value = "'23"
$ws.Range("G32").Value = "'23"
$ws.Range("G33").Value = "'23"
$ws.Range("G34").Value = "'23"
$ws.Range("G35").Value = "'23"
$ws.Range("G36").Value = "'23"
$ws.Range("G37").Value = "'23"
$ws.Range("G38").Value = "'23"
$ws.Range("G39").Value = "'23"
$ws.Range("D40").Value = "'0.04698"
$ws.Range("G40").Value = "'23"
$ws.Range("G41").Value = "'23"
$ws.Range("D42").Value = "'0.1161"
$ws.Range("G42").Value = "'23"
$ws.Range("D43").Value = "'0.003351"
$ws.Range("G43").Value = "'23"
$ws.Range("D44").Value = "'0.01168"
$ws.Range("G44").Value = "'23"
$ws.Range("D45").Value = "'0.00006285"
$ws.Range("G45").Value = "'23"
$ws.Range("D46").Value = "'0.0009897"
$ws.Range("G46").Value = "'23"
$ws.Range("G47").Value = "'23"
$ws.Range("D48").Value = "'0.9197"
$ws.Range("G48").Value = "'23"
$ws.Range("D49").Value = "'0.002107"
$ws.Range("G49").Value = "'23"
$ws.Range("D50").Value = "'0.00001399"
$ws.Range("G50").Value = "'23"
$ws.Range("D51").Value = "'0.01239"
$ws.Range("G51").Value = "'23"
